$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in Milestone (E) and Completed(X) (F) columns for the rows that were
# previously blank. Column E gets "II" or "III" depending on which milestone
# the feature applies to; column F always gets "X".
$milestoneII = @(5, 6, 7, 8, 9, 10, 11, 14, 40)
$milestoneIII = @(15, 24, 30, 82)

foreach ($r in $milestoneII) {
    $ws.Range("E$r").Value = "II"
    $ws.Range("F$r").Value = "X"
}

foreach ($r in $milestoneIII) {
    $ws.Range("E$r").Value = "III"
    $ws.Range("F$r").Value = "X"
}

# Mark the extra-credit rows (Effective Use of GIT / Student Confidence) as
# completed for milestones II and III as well.
$ws.Range("D91").Value = "X"
$ws.Range("E91").Value = "X"
$ws.Range("D92").Value = "X"
$ws.Range("E92").Value = "X"

# Update the sheet view: scroll back to the top-left and move the selection.
$ws.Activate()
[void]$ws.Range("A1").Select()
[void]$ws.Range("E22").Select()
